# Generate Report for Handoff
# Updates the localization-status report after a new handoff run:
#  - Priority column ("ht") is now set for the rows that were handed off
#  - The "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" columns
#    are refreshed to the new handoff timestamp for those same rows.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows 7-12 on each sheet correspond to the six files that were just
# (re)handed off.
$wsOverview.Range("G7:G12").Value = "2016-08-17 02:20:20"

$wsZhCn.Range("E7:E12").Value = "ht"
$wsZhCn.Range("H7:H12").Value = "2016-08-17 02:20:14"

$wsDeDe.Range("E7:E12").Value = "ht"
$wsDeDe.Range("H7:H12").Value = "2016-08-17 02:20:20"
